$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entity")

# Rename field "Id" -> "RefId" (column D, row 3 of the Entity field-definition table)
$ws.Range("D3").Value = "RefId"

# Update the cell selection to match the authored state (cosmetic)
$ws.Range("G12").Select()
